# QA completion: triplicate the last (row 18) record of patient Teresa
# Asensio Navarro into rows 19-21, each stamped with its own submission
# timestamp in column BT ("timestamp"), matching the three extra QA
# register entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common field values shared by rows 19-21 (identical to row 18's record).
# Kind drives how the value is written:
#   "str"  -> text column (written via a NumberFormat detour so it keeps
#             its General-number-format style instead of minting one)
#   "num"  -> plain numeric column
#   "date" -> numeric date-serial column, formatted like the existing
#             date columns (AD/AL/BF) so it reuses their style
$rowData = @(
    @{ Col="Q"; Kind="str"; Val='Segmentectomia o Bisegmentectomia' },
    @{ Col="R"; Kind="str"; Val='bisegmentectomia 6/7' },
    @{ Col="S"; Kind="str"; Val='Bisegmentectomia6i7' },
    @{ Col="T"; Kind="num"; Val='1703' },
    @{ Col="U"; Kind="str"; Val='13/12/2019' },
    @{ Col="V"; Kind="str"; Val='13/12/2019' },
    @{ Col="W"; Kind="str"; Val='10/12/2019' },
    @{ Col="X"; Kind="str"; Val='VIII,V' },
    @{ Col="Y"; Kind="num"; Val='1555' },
    @{ Col="Z"; Kind="str"; Val='Teresa' },
    @{ Col="AA"; Kind="str"; Val='Asensio' },
    @{ Col="AB"; Kind="str"; Val='Navarro' },
    @{ Col="AC"; Kind="num"; Val='11396316' },
    @{ Col="AD"; Kind="date"; Val='43900' },
    @{ Col="AE"; Kind="str"; Val='Si' },
    @{ Col="AF"; Kind="str"; Val='Dona' },
    @{ Col="AG"; Kind="str"; Val='68' },
    @{ Col="AH"; Kind="str"; Val='68' },
    @{ Col="AI"; Kind="num"; Val='159' },
    @{ Col="AJ"; Kind="num"; Val='27' },
    @{ Col="AK"; Kind="num"; Val='2' },
    @{ Col="AL"; Kind="date"; Val='43747' },
    @{ Col="AM"; Kind="str"; Val='No' },
    @{ Col="AN"; Kind="str"; Val='Resecció Menor (<3 segm)' },
    @{ Col="AO"; Kind="str"; Val='1er temps (mobilització)' },
    @{ Col="AP"; Kind="str"; Val='No' },
    @{ Col="AQ"; Kind="str"; Val='No' },
    @{ Col="AR"; Kind="str"; Val='Impressió R0' },
    @{ Col="AS"; Kind="num"; Val='1' },
    @{ Col="AT"; Kind="num"; Val='1.8' },
    @{ Col="AU"; Kind="str"; Val='No' },
    @{ Col="AV"; Kind="str"; Val='No' },
    @{ Col="AW"; Kind="str"; Val='No' },
    @{ Col="AX"; Kind="str"; Val='No' },
    @{ Col="AY"; Kind="str"; Val='No' },
    @{ Col="AZ"; Kind="str"; Val='0' },
    @{ Col="BA"; Kind="num"; Val='0' },
    @{ Col="BB"; Kind="num"; Val='1' },
    @{ Col="BC"; Kind="num"; Val='1' },
    @{ Col="BD"; Kind="num"; Val='0.3' },
    @{ Col="BE"; Kind="str"; Val='No' },
    @{ Col="BF"; Kind="date"; Val='44557' },
    @{ Col="BG"; Kind="str"; Val='No' },
    @{ Col="BH"; Kind="str"; Val='No' },
    @{ Col="BI"; Kind="str"; Val='Viu' },
    @{ Col="BJ"; Kind="str"; Val='No' },
    @{ Col="BK"; Kind="str"; Val='No' },
    @{ Col="BL"; Kind="str"; Val='No' },
    @{ Col="BM"; Kind="str"; Val='No' },
    @{ Col="BN"; Kind="str"; Val='No' },
    @{ Col="BO"; Kind="str"; Val='No' },
    @{ Col="BP"; Kind="str"; Val='No' },
    @{ Col="BQ"; Kind="str"; Val='No' },
    @{ Col="BR"; Kind="str"; Val='No' },
    @{ Col="BS"; Kind="num"; Val='6' }
)

# Per-row submission timestamps for column BT ("timestamp").
$timestamps = @(44611.78876434028, 44611.79149173611, 44611.79204508102)

# Columns that hold numeric-looking text (e.g. "68", "0") need a quick
# NumberFormat detour so Excel stores them as text instead of silently
# coercing to a number. Every cell then gets its NumberFormat (re)applied
# explicitly (General for text/number columns, m/d/yy for date columns)
# so it is assigned a concrete style matching row 18's (General-numFmt /
# date-numFmt) cell styles instead of being left style-less.
$textCells = @()

for ($i = 0; $i -lt 3; $i++) {
    $r = 19 + $i

    foreach ($field in $rowData) {
        $cell = $ws.Range($field.Col + $r)
        if ($field.Kind -eq "date") {
            $cell.Value = [double]$field.Val
            $cell.NumberFormat = "m/d/yy"
        } elseif ($field.Kind -eq "num") {
            $cell.Value = [double]$field.Val
            $cell.NumberFormat = "General"
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = $field.Val
            $textCells += $cell
        }
    }

    $bt = $ws.Range("BT" + $r)
    $bt.Value = $timestamps[$i]
    $bt.NumberFormat = "m/d/yy"
}

foreach ($cell in $textCells) {
    $cell.NumberFormat = "General"
}
